# Correct worksheet name in template: "Data Quality" -> "Provider Submissions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Provider Submissions"
